# Updated cryptos list with latest Price (D) and Volume(1h) (E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is a plain decimal number (e.g. "0.9969").
# Excel would otherwise auto-convert these strings to numeric values, so
# we briefly mark them as Text, write the literal string, then restore the
# default "Normal" style so formatting matches the original sheet.
$plainNumberPriceCells = @("D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D17","D19","D21","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($c in $plainNumberPriceCells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "25.602.66"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.670.02"
$ws.Range("E3").Value = "  +0.98%  "
$ws.Range("D4").Value = "0.9969"
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").Value = "237.84"
$ws.Range("E5").Value = "  +0.78%  "
$ws.Range("D6").Value = "0.9979"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "0.4823"
$ws.Range("D8").Value = "0.2646"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "0.06172"
$ws.Range("E9").Value = "  +3.16%  "
$ws.Range("D10").Value = "0.07107"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").Value = "1.664.32"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "14.96"
$ws.Range("E12").Value = "  +4.32%  "
$ws.Range("D13").Value = "0.6013"
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").Value = "4.439"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").Value = "74.73"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "0.9976"
$ws.Range("E17").Value = "  -0.57%  "
$ws.Range("D18").Value = "25.558.08"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "0.000006811"
$ws.Range("E19").Value = "  +4.41%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "4.484"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").Value = "1.877.23"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("D23").Value = "8.722"
$ws.Range("E23").Value = "  +2.89%  "
$ws.Range("D24").Value = "5.400"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").Value = "134.53"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "15.12"
$ws.Range("E26").Value = "  +2.77%  "
$ws.Range("D27").Value = "1.404"
$ws.Range("E27").Value = "  +1.69%  "
$ws.Range("D28").Value = "104.86"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("D29").Value = "1.712"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("D30").Value = "3.969"
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("D31").Value = "3.701"
$ws.Range("E31").Value = "  +5.06%  "
$ws.Range("D32").Value = "0.07683"
$ws.Range("E32").Value = "  -2.52%  "
$ws.Range("D33").Value = "0.04367"
$ws.Range("E33").Value = "  -5.10%  "
$ws.Range("D34").Value = "0.9972"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "2.618"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").Value = "0.6223"
$ws.Range("E36").Value = "  +6.62%  "
$ws.Range("D37").Value = "0.9533"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").Value = "2.622"
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("D39").Value = "0.8644"
$ws.Range("E39").Value = "  +2.66%  "
$ws.Range("D40").Value = "0.9979"
$ws.Range("E40").Value = "  -0.54%  "
$ws.Range("D41").Value = "0.01515"
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").Value = "1.875"
$ws.Range("E42").Value = "  +2.02%  "
$ws.Range("D43").Value = "98.37"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "0.3798"
$ws.Range("E44").Value = "  +2.79%  "
$ws.Range("D45").Value = "4.723"
$ws.Range("E45").Value = "  -2.53%  "
$ws.Range("D46").Value = "0.1124"
$ws.Range("E46").Value = "  -0.85%  "
$ws.Range("D47").Value = "6.257"
$ws.Range("E47").Value = "  +3.62%  "
$ws.Range("D48").Value = "0.05259"
$ws.Range("E48").Value = "  +2.04%  "
$ws.Range("D49").Value = "29.65"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "7.402"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("D51").Value = "0.3365"
$ws.Range("E51").Value = "  +1.79%  "

# Restore default styling on the cells we temporarily formatted as Text.
foreach ($c in $plainNumberPriceCells) { $ws.Range($c).Style = "Normal" }
